# Error Calculations and Plots
# Apply the recorded edits to the missing_data worksheet:
#   1. Remove the "RM 232" (row 26) and "SC 92" (row 28) samples entirely
#      (rows shift up for everything below them).
#   2. Fill in / clear a handful of individual C/D/E cells to reflect the
#      corrected (re-imputed) values for this combination/seed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two removed-sample rows -------------------------------
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("28:28").Delete()   # "SC 92"
$ws.Rows("26:26").Delete()   # "RM 232"

# --- 2. Per-cell value corrections (row numbers below are POST-delete) ---
$ws.Range("D2").Value = -13.5
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = -6.4
$ws.Range("E5").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("D12").Value = -14.1
$ws.Range("D14").Value = ""
$ws.Range("E15").Value = -8.4
$ws.Range("E18").Value = -8.5
$ws.Range("E19").Value = ""
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("E22").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("E23").Value = -7
$ws.Range("D24").Value = ""
$ws.Range("E25").Value = -7.1
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("E27").Value = ""
$ws.Range("C30").Value = 11.4
$ws.Range("D31").Value = -13.7
$ws.Range("C32").Value = ""
$ws.Range("D33").Value = -14.1
